$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at row 7. This pushes the two blank spacer rows (old
#    7 & 8) and the merged note row (old 9) down by one, and the new row 7
#    inherits the formatting of row 6 (matches native Excel row-insert
#    behaviour).
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).Insert()

# Copy the (inherited) formatting of the brand-new row 7 onto row 8 as well,
# since row 8 will also become a populated data row (old row 7/8 blank
# spacer rows collapse into a single blank row at row 9).
$ws.Range("A7:L7").Copy()
$ws.Range("A8:L8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Row 4 - J1 connector (plug). Reference changes from "J1, J2" to just
#    "J1"; quantity drops from 2 to 1 (L4 formula recalculates automatically).
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "J1"
$ws.Range("J4").Value = 1

# ---------------------------------------------------------------------------
# 3. Row 5 - J2 connector. Previously this row held the J3* header; it now
#    becomes a brand-new line item for J2, using a Molex receptacle part
#    instead of the plug used for J1.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "J2"
$ws.Range("C5").Value = "CON-71439-2164"
$ws.Range("D5").Value = "CON-71439-2164"
$ws.Range("E5").Value = "Molex"
$ws.Range("F5").Value = "71439-2164"
$ws.Range("G5").Value = "Digi-Key"
$ws.Range("G5").Style = $ws.Range("G4").Style
$ws.Range("H5").Value = "WM17222-ND"
$ws.Range("I5").Value = "CONN RECPT 64POS VERT 1MM SMD"
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 5.74

# ---------------------------------------------------------------------------
# 4. Row 6 - J3* header (the part that used to live in row 5, now shifted
#    down one row because of the new J2 line item above it).
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "J3*"
$ws.Range("B6").Style = $ws.Range("B4").Style
$ws.Range("C6").Value = "CON-2X8"
$ws.Range("D6").Value = "CON-2X8"
$ws.Range("E6").Value = "Samtec"
$ws.Range("F6").Value = "TSW-150-08-T-D-RA"
$ws.Range("G6").Value = "Digi-Key"
$ws.Range("G6").Style = $ws.Range("G5").Style
$ws.Range("H6").Value = "SAM1049-50-ND"
$ws.Range("I6").Value = "CONN HEADR 100PS .100 DL R/A TIN"
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 4.62

# ---------------------------------------------------------------------------
# 5. Row 7 (new) - Test points TP1, TP2 (TP11 removed from the reference
#    designator list).
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "TP1, TP2"
$ws.Range("C7").Value = "TESTPOINT"
$ws.Range("D7").Value = "TESTPOINT-YELLOW"
$ws.Range("E7").Value = "Keystone"
$ws.Range("F7").Value = 5009
$ws.Range("G7").Value = "Digi-Key"
$ws.Range("H7").Value = "5009K-ND"
$ws.Range("I7").Value = "TEST POINT PC COMPACT .063""D YLW"
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 0.36
$ws.Range("L7").Formula = "=J7*K7"

# ---------------------------------------------------------------------------
# 6. Row 8 (new) - Standoffs M1-M4.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "M1,M2,M3,M4"
$ws.Range("C8").Value = "STANDOFF"
$ws.Range("D8").Value = "Aluminum"
$ws.Range("E8").Value = "Keystone"
$ws.Range("F8").Value = 8401
$ws.Range("G8").Value = "Digi-Key"
$ws.Range("H8").Value = "8401K-ND"
$ws.Range("I8").Value = "STDOFF HEX M/F 4-40 .500""L ALUM"
$ws.Range("J8").Value = 10
$ws.Range("K8").Value = 0.717
$ws.Range("L8").Formula = "=J8*K8"

# ---------------------------------------------------------------------------
# 7. Cosmetic updates: column B got wider to fit "M1,M2,M3,M4", and the
#    selection cursor ended up on E17.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 19.42578125
$ws.Range("E17").Select()
